# Lesson 9.1 Classes, Objects, and Methods — "starting to work on Module 09"
#
# 1) Refresh the cached date-placeholder text (10/29/2016 -> 11/1/2016) on
#    every slide layout, the slide master, and the notes master.
# 2) Slide 18 ("foo ... is a function of no arguments ..."): reword the
#    explanatory rectangle and resize/reposition it.
# 3) Slide 21 ("Here's the definition of Class2% ..."): bold the class name
#    "Class2% " inside the sentence.

$p = $ppt.ActivePresentation

$oldDate = "10/29/2016"
$newDate = "11/1/2016"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# -- Slide master's own date placeholder --
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# -- Every slide layout's date placeholder --
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# -- Notes master's date placeholder --
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# -- Slide 18: reword + resize the "foo" explanation rectangle --
$s18 = $p.Slides.Item(18)
$rect18 = $s18.Shapes.Item(10)

$rect18.Left = 176.1082
$rect18.Top = 481.22551181102364
$rect18.Width = 367.56535433070866
$rect18.Height = 51.8775

$tr18 = $rect18.TextFrame.TextRange
$middle18 = $tr18.Characters(4, 42)
$middle18.Text = " will be the name of a method that takes no arguments (legal in #"

# -- Slide 21: bold "Class2% " within the sentence --
$s21 = $p.Slides.Item(21)
$rect21 = $s21.Shapes.Item(4)
$tr21 = $rect21.TextFrame.TextRange
$className = $tr21.Characters(26, 8)
$className.Font.Bold = $true
